$d = $word.ActiveDocument

# The final paragraph of the scene ("However, I only manage to get a few
# scribbles down by the end of class .") is currently split across five
# separate runs (an artifact of earlier word-by-word edits). Collapse it
# back down into a single run while leaving the run formatting (sz 24,
# szCs 24, rtl 0) and the surrounding paragraph untouched.

$old = "However, I only manage to get a few scribbles down by the end of class ."
$new = "However, I only manage to get a few scribbles down by the end of class ."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "done"
